# Updated diagrams in DG
#
# This slide is the "UiComponentClassDiagram" last-slide variant (creationId
# 2776882492). Two labelled rectangles get renamed/resized (Person* ->
# ModuleTaken*) and four connectors attached to them are re-routed to match
# the new geometry.
#
# NOTE on the numeric literals below: PowerPoint's Shape.Left/Top/Width/Height
# are 32-bit (single-precision) floats expressed in points. To land on an
# exact target EMU value after the point->EMU round-trip (EMU = points/72*914400,
# truncated) we pick the double literal whose nearest float32 converts back to
# the desired EMU value rather than just emu/914400*72 (which can be off by
# one EMU because of the float32 truncation).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($sh in $s.Shapes) {

    if ($sh.Id -eq 36) {
        # "PersonListPanel" rectangle -> becomes a two-line "ModuleTaken" / "ListPanel" label.
        # Position/size change: off (2592526,3991960)->(2592526,3962400) ext (1093635,236841)->(1093635,266401)
        $sh.Top = 312.0000305175781
        $sh.Height = 20.976457595825195

        $tr = $sh.TextFrame.TextRange
        $tr.Text = "ModuleTaken"
        $tr.InsertAfter([char]13 + "ListPanel")
    }

    if ($sh.Id -eq 37) {
        # "PersonCard" rectangle -> "ModuleTakenCard"
        # Position/size change: off (3839323,4228801)->(3839322,4228801) ext (1040906,236841)->(1266077,236841)
        $sh.Left = 302.308837890625
        $sh.Width = 99.69110870361328

        $sh.TextFrame.TextRange.Text = "ModuleTakenCard"
    }

    if ($sh.Id -eq 47) {
        # Elbow Connector 63 (stCxn 39 -> endCxn 36), re-routed after shape 36 moved.
        # off (1883148,3401003)->(1890538,3393613) ext (1242356,176400)->(1227576,176400)
        $sh.Left = 148.86126708984375
        $sh.Top = 267.213623046875
        $sh.Width = 96.65953063964844
    }

    if ($sh.Id -eq 77) {
        # Elbow Connector 63 (stCxn 16 -> endCxn 37), re-routed after shape 37 resized.
        # off (4174488,2991741)->(4287073,3104326) ext (2061222,649740)->(2061222,424570)
        $sh.Left = 337.5648193359375
        $sh.Top = 244.43511962890625
        $sh.Height = 33.43070983886719
    }

    if ($sh.Id -eq 137) {
        # Elbow Connector 136 (stCxn 36 -> endCxn 37), tiny 1-EMU geometry nudge.
        # off (3430123,3938021)->(3430123,3938022) ext (118421,699979)->(118421,699978)
        $sh.Top = 310.0804748535156
        $sh.Height = 55.11637878417969
    }

    if ($sh.Id -eq 140) {
        # Elbow Connector 63 (stCxn 16 -> endCxn 36), re-routed after shape 36 moved.
        # off (3695875,2276286)->(3703265,2268896) ext (1824381,1843808)->(1809601,1843808)
        $sh.Left = 291.5956726074219
        $sh.Top = 178.65322875976562
        $sh.Width = 142.48828125
    }
}
